# Add a "Score Type" column to the "Clinical Data" sheet of the
# MDS-UPDRS III template, right after the PatientID column (i.e. as the
# new column B), shifting every existing column one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clinical Data")

# Insert a new, blank column before the current column B ("Timeline").
# This shifts every column from B..AJ one place to the right (C..AK) and
# carries the header row's formatting along with it.
$ws.Columns.Item(2).EntireColumn.Insert()

# Give the new column its header text.
$ws.Cells.Item(1, 2).Value = "Score Type"

# The worksheet used to have one extra, already-blank formatted column
# past the last real header (originally column AJ, with no content) that
# only existed to carry column-level formatting. After the insert above
# it now sits at column 37 (AK). Drop its formatting/content so the
# sheet's used range again ends exactly at the last real header column.
$ws.Columns.Item(37).ClearFormats()
$ws.Cells.Item(1, 37).ClearContents()

# Match the saved selection/view state.
$ws.Range("C7").Select()
